# feat: add 2022-Q4 data
#
# 1) "总计" sheet: insert a new row (2022-Q4 totals) above the existing
#    2022-Q2 totals row.
# 2) Worksheets: insert a brand-new "2022-Q4" sheet (fund-level detail)
#    positioned between "总计" and "2022-Q2"; "2022-Q2" keeps its
#    original data (rebuilt in-place so sheet order/ids line up).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert the new 2022-Q4 row above row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows(2).Insert()

# Pull formatting (bold/border/centered) for the new A2 from A3 (the row
# we just pushed down, which already carries the correct style).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.06

# The old row (now row 3) was rank/index 0 in the "总计" table; since a
# newer quarter was inserted above it, its index becomes 1.
$summary.Range("A3").Value = 1

# ---------------------------------------------------------------------
# Step 2: fund-detail sheets.
#
# The sheet-id bookkeeping only lines up with the target layout
# (总计=1, 2022-Q4=2, 2022-Q2=3) if "2022-Q2" is removed and both
# replacement sheets are (re)created fresh and in order - any
# Worksheet.Copy() call burns an id even if that sheet is later
# deleted, so the original "2022-Q2" content is re-entered verbatim
# instead of copied.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$oldQ2.Delete()

$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q4.Name = "2022-Q4"

$q2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$q2.Name = "2022-Q2"

# --- Populate "2022-Q4" ------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$q4.Range("B1:H1").Font.Bold = $true
$q4.Range("B1:H1").Borders.LineStyle = 1
$q4.Range("B1:H1").HorizontalAlignment = -4108
$q4.Range("B1:H1").VerticalAlignment = -4160

$q4.Range("A2").Value = 0
$q4.Range("A2").Font.Bold = $true
$q4.Range("A2").Borders.LineStyle = 1
$q4.Range("A2").HorizontalAlignment = -4108
$q4.Range("A2").VerticalAlignment = -4160

$q4.Range("B2").Value = "'233009"
$q4.Range("C2").Value = "大摩多因子精选策略混合"
$q4.Range("D2").Value = "'6.42"
$q4.Range("E2").Value = "'91.11"
$q4.Range("F2").Value = "'0.98"
$q4.Range("G2").Value = "'0.0629"
$q4.Range("H2").Value = 4

# --- Populate "2022-Q2" (verbatim re-entry of the pre-existing data) ---
$q2.Range("B1").Value = "基金代码"
$q2.Range("C1").Value = "基金名称"
$q2.Range("D1").Value = "基金规模"
$q2.Range("E1").Value = "股票总仓位"
$q2.Range("F1").Value = "仓位占比"
$q2.Range("G1").Value = "持有市值(亿元)"
$q2.Range("H1").Value = "仓位排名"
$q2.Range("B1:H1").Font.Bold = $true
$q2.Range("B1:H1").Borders.LineStyle = 1
$q2.Range("B1:H1").HorizontalAlignment = -4108
$q2.Range("B1:H1").VerticalAlignment = -4160

$q2.Range("A2").Value = 0
$q2.Range("A2").Font.Bold = $true
$q2.Range("A2").Borders.LineStyle = 1
$q2.Range("A2").HorizontalAlignment = -4108
$q2.Range("A2").VerticalAlignment = -4160

$q2.Range("B2").Value = "'002871"
$q2.Range("C2").Value = "华夏智胜价值成长股票A"
$q2.Range("D2").Value = "'2.60"
$q2.Range("E2").Value = "'93.26"
$q2.Range("F2").Value = "'0.96"
$q2.Range("G2").Value = "'0.0250"
$q2.Range("H2").Value = 10

$q2.Range("A3").Value = 1
$q2.Range("A3").Font.Bold = $true
$q2.Range("A3").Borders.LineStyle = 1
$q2.Range("A3").HorizontalAlignment = -4108
$q2.Range("A3").VerticalAlignment = -4160

$q2.Range("B3").Value = "'002872"
$q2.Range("C3").Value = "华夏智胜价值成长股票C"
$q2.Range("D3").Value = "'1.69"
$q2.Range("E3").Value = "'93.26"
$q2.Range("F3").Value = "'0.96"
$q2.Range("G3").Value = "'0.0162"
$q2.Range("H3").Value = 10

Write-Output "done"
